$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) — update F5, F6, F7
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 507
$ws1.Range("F6").Value = 1509
$ws1.Range("F7").Value = 1001

# Sheet "全部类型" (sheet4) — update F5, F6, F8
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 507
$ws4.Range("F6").Value = 1509
$ws4.Range("F8").Value = 1001
